{"js": "// Add new \"Laptop\" sub-category products (display light / elegant laptops)\n// right after the existing \"Acer: 1\" entry (end of \"H\u1ecdc t\u1eadp, v\u0103n ph\u00f2ng\" group)\n// and before the next top-level \"Table: 1\" item.\n//\n// New structure appended (mirrors existing siblings: category at list level 6,\n// brand/qty entries at list level 7, both using style \"ListParagraph\" and\n// the same numbering definition, numId 4):\n//   \u0110\u1ed3 ho\u1ea1, k\u1ef9 thu\u1eadt   (level 6)\n//     Hp: 1            (level 7)\n//   M\u1ecfng nh\u1eb9           (level 6)\n//     Lenovo: 1        (level 7)\n//   Sang tr\u1ecdng         (level 6)\n//     Lenovo: 1        (level 7)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the anchor paragraph: \"Acer: 1\" immediately followed by \"Table: 1\".\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length - 1; i++) {\n  if (\n    paragraphs.items[i].text === \"Acer: 1\" &&\n    paragraphs.items[i + 1].text === \"Table: 1\"\n  ) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not find anchor paragraph \"Acer: 1\" before \"Table: 1\".');\n}\n\n// New items to insert, in order, each with its target outline level.\nconst newItems = [\n  { text: \"\u0110\u1ed3 ho\u1ea1, k\u1ef9 thu\u1eadt\", level: 6 },\n  { text: \"Hp: 1\", level: 7 },\n  { text: \"M\u1ecfng nh\u1eb9\", level: 6 },\n  { text: \"Lenovo: 1\", level: 7 },\n  { text: \"Sang tr\u1ecdng\", level: 6 },\n  { text: \"Lenovo: 1\", level: 7 }\n]; // NB: brand/qty runs mirror existing sibling items verbatim\n\nlet previous = anchor;\nfor (const item of newItems) {\n  const newPara = previous.insertParagraph(item.text, Word.InsertLocation.after);\n  newPara.load(\"listItemOrNullObject\");\n  await context.sync();\n\n  newPara.listItemOrNullObject.level = item.level;\n  await context.sync();\n\n  previous = newPara;\n}\n", "ps1": "# Add new \"Laptop\" sub-category products (display light / elegant laptops)\n# right after the existing \"Acer: 1\" entry (end of the \"Hoc tap, van phong\"\n# group) and before the next top-level \"Table: 1\" item.\n#\n# New structure appended (mirrors existing siblings: category at list level\n# 7 i.e. ilvl=6, brand/qty entries at list level 8 i.e. ilvl=7 - both 1-based\n# ListLevelNumber - using style \"ListParagraph\" and the same numbering\n# definition, numId 4):\n#   Do hoa, ky thuat   (level 7 / ilvl 6)  -- \"Do hoa\" = graphics/design\n#     Hp: 1            (level 8 / ilvl 7)\n#   Mong nhe           (level 7 / ilvl 6)\n#     Lenovo: 1        (level 8 / ilvl 7)\n#   Sang trong         (level 7 / ilvl 6)\n#     Lenovo: 1        (level 8 / ilvl 7)\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph: \"Acer: 1\" immediately followed by \"Table: 1\".\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -eq \"Acer: 1`r\") {\n    $next = $d.Paragraphs.Item($i + 1)\n    if ($next.Range.Text -eq \"Table: 1`r\") {\n      $anchorIndex = $i\n      break\n    }\n  }\n}\n\nif ($anchorIndex -eq -1) {\n  throw 'Could not find anchor paragraph \"Acer: 1\" before \"Table: 1\".'\n}\n\n# New items to insert, in order, each with its target 1-based list level.\n$items = @(\n  @{ Text = \"\u0110\u1ed3 ho\u1ea1, k\u1ef9 thu\u1eadt\"; Level = 7 },\n  @{ Text = \"Hp: 1\"; Level = 8 },\n  @{ Text = \"M\u1ecfng nh\u1eb9\"; Level = 7 },\n  @{ Text = \"Lenovo: 1\"; Level = 8 },\n  @{ Text = \"Sang tr\u1ecdng\"; Level = 7 },\n  @{ Text = \"Lenovo: 1\"; Level = 8 }\n)\n\n$insertAt = $anchorIndex\nforeach ($item in $items) {\n  $p = $d.Paragraphs.Item($insertAt)\n  $p.Range.InsertParagraphAfter()\n  $insertAt = $insertAt + 1\n\n  $newPara = $d.Paragraphs.Item($insertAt)\n  $newPara.Range.Text = $item.Text\n  $newPara.Range.ListFormat.ListLevelNumber = $item.Level\n}\n"}
